# Update Name of Algo
# Apply updated KNN-imputed values in column E for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 16.272
    14 = 16.95500000000001
    21 = 16.541
    23 = 16.55
    25 = 17.493
    26 = 16.42
    29 = 16.85
    53 = 16.639
    57 = 16.556
    59 = 16.305
    69 = 17.563
    79 = 17.06
    83 = 16.535
    91 = 17.612
    93 = 17.291
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
